# Update script for "Fruta, Femacal de La Calera - Pera" sheet
# Applies the weekly roll-forward described in the commit diff ("Fruta /
# hortaliza, semanal"):
#   - rows 775-806: date / variety / quality / volume / price columns
#     shifted down to reflect the next week's report
#   - rows 807-810: newly appended rows for the latest week (date 44491:
#     two extra Packham's Triumph grades plus the Winter Nelis entries
#     that used to sit in row 806)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=775; D=44509; M=78; N=13000; O=13000; P=13000; S=722},
    @{Row=776; D=44509; M=50; N=14000; O=14000; P=14000; S=778},
    @{Row=777; D=44509; M=67; N=12000; O=12000; P=12000; S=667},
    @{Row=778; D=44509; M=68; N=11000; O=11000; P=11000; S=611},
    @{Row=779; D=44421; N=11000; O=11000; P=11000; S=611},
    @{Row=780; D=44421; M=70; N=12000; O=12000; P=12000; S=667},
    @{Row=781; D=44421; N=10000; O=10000; P=10000; S=556},
    @{Row=782; D=44421; M=70; N=9000; O=9000; P=9000; S=500},
    @{Row=783; D=44383; K='Packham''s Triumph'; L='Especial'; M=80; N=10000; O=10000; P=10000; S=556},
    @{Row=784; D=44383; K='Packham''s Triumph'; L='Extra (doble especial)'; M=68; N=11000; O=11000; P=11000; S=611},
    @{Row=785; D=44383; L='Primera'; M=80; N=9000; O=9000; P=9000; S=500},
    @{Row=786; D=44383; L='Segunda'; M=75; N=8000; O=8000; P=8000; S=444},
    @{Row=787; D=44244; K='Bartlett de verano'; M=75; N=12000; O=12000; P=12000; S=667},
    @{Row=788; D=44244; K='Bartlett de verano'; M=50; N=10000; O=10000; P=10000; S=556},
    @{Row=789; D=44307; M=85; N=12000; O=12000; P=12000; S=667},
    @{Row=790; D=44307; L='Extra (doble especial)'; M=50; N=13000; O=13000; P=13000; S=722},
    @{Row=791; D=44307; L='Primera'; M=87; N=10000; O=10000; P=10000; S=556},
    @{Row=792; D=44307; L='Segunda'; M=80; N=9000; O=9000; P=9000; S=500},
    @{Row=793; D=44273; L='Especial'; N=13000; O=13000; P=13000; S=722},
    @{Row=794; D=44273; N=11000; O=11000; P=11000; S=611},
    @{Row=795; D=44273},
    @{Row=796; D=44433; M=78; N=11000; O=11000; P=11000; S=611},
    @{Row=797; D=44433; M=70; N=12000; O=12000; P=12000; S=667},
    @{Row=798; D=44433; M=75},
    @{Row=799; D=44433; L='Segunda'; M=70; N=9000; O=9000; P=9000; S=500},
    @{Row=800; D=44302; L='Especial'; M=65; N=12000; O=12000; P=12000; S=667},
    @{Row=801; D=44302; L='Extra (doble especial)'; M=68},
    @{Row=802; D=44302; L='Primera'; M=70; N=10000; O=10000; P=10000; S=556},
    @{Row=803; D=44179; L='Especial'; N=15000; O=15000; P=15000; S=833},
    @{Row=804; D=44179; L='Primera'; M=75; N=14000; O=14000; P=14000; S=778},
    @{Row=805; D=44179; K='Packham''s Triumph'; L='Segunda'; M=70; N=13000; O=13000; P=13000; S=722},
    @{Row=806; K='Packham''s Triumph'; L='Especial'; M=75; N=13000; O=13000; P=13000; S=722}
)

foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @('D','K','L','M','N','O','P','S')) {
        if ($u.ContainsKey($col)) {
            $ws.Range("$col$r").Value = $u[$col]
        }
    }
}

# Shared (constant-across-rows) column values for this market/product block
$const = @{
    A = 3
    B = 'Femacal de La Calera'
    C = 'Coquimbo'
    E = 5
    F = 'Fruta'
    G = 100104
    H = 'Frutos de pepita'
    I = 100104005
    J = 'Pera'
    Q = '$/caja 18 kilos empedrada'
    R = 'Región de O''Higgins'
    T = 18
}

$newRows = @(
    @{Row=807; D=44491; K='Packham''s Triumph'; L='Primera'; M=78; N=12000; O=12000; P=12000; S=667},
    @{Row=808; D=44491; K='Packham''s Triumph'; L='Segunda'; M=70; N=1000; O=1000; P=1000; S=56},
    @{Row=809; D=44491; K='Winter Nelis'; L='Primera'; M=50; N=12000; O=12000; P=12000; S=667},
    @{Row=810; D=44491; K='Winter Nelis'; L='Segunda'; M=58; N=10000; O=10000; P=10000; S=556}
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    foreach ($col in $const.Keys) {
        $ws.Range("$col$r").Value = $const[$col]
    }
    foreach ($col in @('D','K','L','M','N','O','P','S')) {
        $ws.Range("$col$r").Value = $nr[$col]
    }
    $ws.Range("D$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
